$d = $word.ActiveDocument

# Locate the sentence containing the word "to" we need to split ("...allow you to display...").
$anchor = $d.Content
$found = $anchor.Find.Execute("you to display", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Anchor text not found"
}

# "you to display" -> the "to" token starts 4 characters in (after "you ").
$toStart = $anchor.Start + 4

# Split range for the "t" of "to" (1 character) - toggling formatting on/off forces
# Word to materialize this character as its own run even though the final
# formatting matches its neighbours, reproducing the run split from the edit.
$tRange = $d.Range($toStart, $toStart + 1)
$tRange.Bold = 1
$tRange.Bold = 0
